# Update build timestamp from "February 03 2026 17.29.55 EST" to
# "February 03 2026 18.05.36 EST" across the workbook.

$wb = $excel.ActiveWorkbook

$newStamp = "February 03 2026 18.05.36 EST"

$wsAbout = $wb.Worksheets.Item("About")
$wsData  = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet ---
$wsAbout.Range("A2").Value = "Version: Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for AMC Coal Mines, Indonesia, M1339, version 'Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" sheet, column S rows 2-12 ---
$newVersionString = "Coal Mine Boundaries and Methane Sources - version 1.0.0 (built on $newStamp)"
for ($row = 2; $row -le 12; $row++) {
    $wsData.Range("S$row").Value = $newVersionString
}
